# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the quarterly fund-holdings table.
# 2) Prepend a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, positioned just before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Borrow header / index-column formatting (and page margins) from the
# "2021-Q3" sheet, which already has the 8-column layout this quarter
# needs.
$template = $wb.Worksheets.Item("2021-Q3")

$template.Range("B1:H1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$template.Range("A2:A17").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.PageSetup.LeftMargin = $template.PageSetup.LeftMargin
$ws.PageSetup.RightMargin = $template.PageSetup.RightMargin
$ws.PageSetup.TopMargin = $template.PageSetup.TopMargin
$ws.PageSetup.BottomMargin = $template.PageSetup.BottomMargin
$ws.PageSetup.HeaderMargin = $template.PageSetup.HeaderMargin
$ws.PageSetup.FooterMargin = $template.PageSetup.FooterMargin

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows. Columns B-G are stored as plain text (matching the other
# quarter sheets), A/H are numeric.
$rows = @(
    @(0,  "010340", "易方达高质量严选三年持有期混合",         "113.51", "91.97", "4.24", "4.8128", 10),
    @(1,  "501054", "东方红睿泽三年定期开放灵活配置混合A",     "109.00", "95.90", "3.34", "3.6406", 9),
    @(2,  "009576", "东方红智远三年持有期混合",                "66.98",  "92.53", "3.53", "2.3644", 7),
    @(3,  "169104", "东方红睿满沪港深灵活配置混合（LOF）",     "48.91",  "92.40", "4.52", "2.2107", 5),
    @(4,  "118001", "易方达亚洲精选股票(QDII)",                "46.85",  "93.54", "3.57", "1.6725", 9),
    @(5,  "006595", "广发港股通优质增长混合",                  "8.53",   "86.63", "5.15", "0.4393", 5),
    @(6,  "010852", "中欧内需成长混合型证券投资基金A",         "5.23",   "91.46", "3.82", "0.1998", 10),
    @(7,  "005620", "中欧品质消费股票A",                       "3.74",   "90.47", "5.29", "0.1978", 10),
    @(8,  "012447", "华夏互联网龙头混合型证券投资基金A",       "3.32",   "83.95", "4.15", "0.1378", 10),
    @(9,  "161132", "易方达科顺定期开放灵活配置混合",          "1.15",   "93.53", "5.58", "0.0642", 9),
    @(10, "007109", "南方沪港深核心优势混合",                  "1.82",   "87.54", "3.50", "0.0637", 4),
    @(11, "005621", "中欧品质消费股票C",                       "1.11",   "90.47", "5.29", "0.0587", 10),
    @(12, "012448", "华夏互联网龙头混合型证券投资基金C",       "1.28",   "83.95", "4.15", "0.0531", 10),
    @(13, "009017", "银华港股通精选股票",                      "0.91",   "86.12", "4.66", "0.0424", 8),
    @(14, "010853", "中欧内需成长混合型证券投资基金C",         "0.67",   "91.46", "3.82", "0.0256", 10),
    @(15, "011032", "东方红睿泽三年定期开放灵活配置混合C",     "0.35",   "95.90", "3.34", "0.0117", 9)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("B$r").ClearFormats()

    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("C$r").ClearFormats()

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("D$r").ClearFormats()

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("E$r").ClearFormats()

    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("F$r").ClearFormats()

    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("G$r").ClearFormats()

    $ws.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to the "总计" sheet.
#    (Re-fetch by name: inserting a sheet before it shifted its tab
#    position, and sheet handles above resolve positionally.)
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows("2:2").Insert()

$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)
$tot.Range("A2").Value = 0

$tot.Range("B2:D2").ClearFormats()
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 16
$tot.Range("D2").Value = 16

# ---------------------------------------------------------------------
# Restore the original active sheet/selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
